# Auto-generated update script
# Commit message: Update automatic: dades i banners [2026-02-24 17:50]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force literal-text writes (via quote-prefix + Copy/PasteSpecial-Values)
# so that numeric- or percent-looking strings (e.g. "40%") are not reinterpreted by Excel
# as numbers, which would otherwise change the cell's style index.
$scratch = $ws.Cells.Item(2, 50)

function Set-TextValue {
    param($TargetCell, [string]$Text)
    $scratch.Value = "'" + $Text
    $scratch.Copy()
    $TargetCell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

Set-TextValue $ws.Range("E2") "2026-02-24 17:48:37"
Set-TextValue $ws.Range("H2") "40%"
Set-TextValue $ws.Range("E3") "2026-02-24 17:48:40"
Set-TextValue $ws.Range("K3") "16.4 MJ/m2"
Set-TextValue $ws.Range("E4") "2026-02-24 17:48:43"
Set-TextValue $ws.Range("J4") "1020.0 hPa"
Set-TextValue $ws.Range("E5") "2026-02-24 17:48:45"
Set-TextValue $ws.Range("K5") "15.7 MJ/m2"
Set-TextValue $ws.Range("E6") "2026-02-24 17:48:48"
Set-TextValue $ws.Range("J6") "1019.9 hPa"
Set-TextValue $ws.Range("O6") "14.2 °C"
Set-TextValue $ws.Range("E7") "2026-02-24 17:48:50"
Set-TextValue $ws.Range("J7") "1020.5 hPa"
Set-TextValue $ws.Range("K7") "15.6 MJ/m2"
Set-TextValue $ws.Range("E8") "2026-02-24 17:48:53"
Set-TextValue $ws.Range("J8") "1019.9 hPa"
Set-TextValue $ws.Range("K8") "16.3 MJ/m2"
Set-TextValue $ws.Range("O8") "16.4 °C"
Set-TextValue $ws.Range("E9") "2026-02-24 17:48:56"
Set-TextValue $ws.Range("O9") "12.1 °C"
Set-TextValue $ws.Range("E10") "2026-02-24 17:48:58"
Set-TextValue $ws.Range("K10") "14.4 MJ/m2"
Set-TextValue $ws.Range("O10") "11.8 °C"
Set-TextValue $ws.Range("E11") "2026-02-24 17:49:01"
Set-TextValue $ws.Range("O11") "9.2 °C"
Set-TextValue $ws.Range("E12") "2026-02-24 17:49:03"
Set-TextValue $ws.Range("E13") "2026-02-24 17:49:05"
Set-TextValue $ws.Range("H13") "63%"
Set-TextValue $ws.Range("J13") "1023.7 hPa"
Set-TextValue $ws.Range("K13") "14.9 MJ/m2"
Set-TextValue $ws.Range("O13") "6.5 °C"
Set-TextValue $ws.Range("E14") "2026-02-24 17:49:08"
Set-TextValue $ws.Range("E15") "2026-02-24 17:49:10"
Set-TextValue $ws.Range("O15") "12.4 °C"
Set-TextValue $ws.Range("E16") "2026-02-24 17:49:13"
Set-TextValue $ws.Range("K16") "13.3 MJ/m2"
Set-TextValue $ws.Range("E17") "2026-02-24 17:49:15"
Set-TextValue $ws.Range("E18") "2026-02-24 17:49:18"
Set-TextValue $ws.Range("H18") "74%"
Set-TextValue $ws.Range("J18") "1020.4 hPa"
Set-TextValue $ws.Range("O18") "11.1 °C"
Set-TextValue $ws.Range("E19") "2026-02-24 17:49:20"
Set-TextValue $ws.Range("O19") "12.9 °C"
Set-TextValue $ws.Range("E20") "2026-02-24 17:49:23"
Set-TextValue $ws.Range("K20") "15.1 MJ/m2"
Set-TextValue $ws.Range("E21") "2026-02-24 17:49:25"
Set-TextValue $ws.Range("J21") "1022.4 hPa"
Set-TextValue $ws.Range("O21") "9.5 °C"
Set-TextValue $ws.Range("E22") "2026-02-24 17:49:28"
Set-TextValue $ws.Range("H22") "20%"
Set-TextValue $ws.Range("K22") "16.4 MJ/m2"
Set-TextValue $ws.Range("E23") "2026-02-24 17:49:30"
Set-TextValue $ws.Range("H23") "22%"
Set-TextValue $ws.Range("K23") "16.2 MJ/m2"
Set-TextValue $ws.Range("O23") "4.9 °C"
Set-TextValue $ws.Range("E24") "2026-02-24 17:49:33"
Set-TextValue $ws.Range("H24") "79%"
Set-TextValue $ws.Range("J24") "1021.7 hPa"
Set-TextValue $ws.Range("K24") "15.0 MJ/m2"
Set-TextValue $ws.Range("O24") "9.8 °C"
Set-TextValue $ws.Range("E25") "2026-02-24 17:49:35"
Set-TextValue $ws.Range("H25") "32%"
Set-TextValue $ws.Range("O25") "7.2 °C"
Set-TextValue $ws.Range("E26") "2026-02-24 17:49:38"
Set-TextValue $ws.Range("K26") "15.7 MJ/m2"
Set-TextValue $ws.Range("E27") "2026-02-24 17:49:40"
Set-TextValue $ws.Range("K27") "15.3 MJ/m2"
Set-TextValue $ws.Range("O27") "6.6 °C"
Set-TextValue $ws.Range("E28") "2026-02-24 17:49:43"
Set-TextValue $ws.Range("J28") "1020.4 hPa"
Set-TextValue $ws.Range("O28") "11.6 °C"
Set-TextValue $ws.Range("E29") "2026-02-24 17:49:46"
Set-TextValue $ws.Range("K29") "14.9 MJ/m2"
Set-TextValue $ws.Range("O29") "10.2 °C"
Set-TextValue $ws.Range("E30") "2026-02-24 17:49:49"
Set-TextValue $ws.Range("E31") "2026-02-24 17:49:51"
Set-TextValue $ws.Range("J31") "1019.5 hPa"
Set-TextValue $ws.Range("E32") "2026-02-24 17:49:54"
Set-TextValue $ws.Range("H32") "65%"
Set-TextValue $ws.Range("O32") "7.8 °C"
Set-TextValue $ws.Range("E33") "2026-02-24 17:49:57"
Set-TextValue $ws.Range("J33") "1022.0 hPa"
Set-TextValue $ws.Range("K33") "15.6 MJ/m2"
Set-TextValue $ws.Range("O33") "8.4 °C"
Set-TextValue $ws.Range("E34") "2026-02-24 17:49:59"
Set-TextValue $ws.Range("H34") "47%"
Set-TextValue $ws.Range("K34") "14.4 MJ/m2"
Set-TextValue $ws.Range("E35") "2026-02-24 17:50:02"
Set-TextValue $ws.Range("J35") "1020.7 hPa"
Set-TextValue $ws.Range("E36") "2026-02-24 17:50:05"
Set-TextValue $ws.Range("J36") "1020.2 hPa"
Set-TextValue $ws.Range("O36") "13.2 °C"
Set-TextValue $ws.Range("E37") "2026-02-24 17:50:07"
Set-TextValue $ws.Range("J37") "1022.4 hPa"
Set-TextValue $ws.Range("O37") "9.0 °C"
Set-TextValue $ws.Range("E38") "2026-02-24 17:50:10"
Set-TextValue $ws.Range("H38") "69%"
Set-TextValue $ws.Range("E39") "2026-02-24 17:50:12"
Set-TextValue $ws.Range("L39") "22.7 km/h - 317º 17:05 TU"
Set-TextValue $ws.Range("N39") "2.3 °C 17:28 TU"
Set-TextValue $ws.Range("O39") "4.8 °C"
Set-TextValue $ws.Range("E40") "2026-02-24 17:50:15"
Set-TextValue $ws.Range("H40") "65%"
Set-TextValue $ws.Range("J40") "1023.1 hPa"
Set-TextValue $ws.Range("O40") "8.7 °C"
Set-TextValue $ws.Range("E41") "2026-02-24 17:50:17"
Set-TextValue $ws.Range("O41") "10.8 °C"
Set-TextValue $ws.Range("E42") "2026-02-24 17:50:20"
Set-TextValue $ws.Range("O42") "11.4 °C"
Set-TextValue $ws.Range("E43") "2026-02-24 17:50:22"
Set-TextValue $ws.Range("H43") "71%"
Set-TextValue $ws.Range("O43") "10.5 °C"
Set-TextValue $ws.Range("E44") "2026-02-24 17:50:25"
Set-TextValue $ws.Range("K44") "15.5 MJ/m2"
Set-TextValue $ws.Range("E45") "2026-02-24 17:50:28"
Set-TextValue $ws.Range("H45") "42%"
Set-TextValue $ws.Range("J45") "1021.0 hPa"
Set-TextValue $ws.Range("L45") "20.2 km/h - 127º 17:20 TU"
Set-TextValue $ws.Range("O45") "10.3 °C"
Set-TextValue $ws.Range("E46") "2026-02-24 17:50:31"
Set-TextValue $ws.Range("H46") "74%"
Set-TextValue $ws.Range("J46") "1021.8 hPa"
Set-TextValue $ws.Range("O46") "10.1 °C"

$scratch.Clear()
